# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that were refreshed when the
# handback report was regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for 024e7ba5-...md (row 3)
$wsOverview.Range("G3").Value = "2016-09-04 20:51:19"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for
# 024e7ba5-... (row 3)
$wsZhCn.Range("H3").Value = "2016-09-04 20:51:14"
$wsZhCn.Range("K3").Value = "2016-09-04 20:51:33"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for
# 024e7ba5-... (row 3)
$wsDeDe.Range("H3").Value = "2016-09-04 20:51:19"
$wsDeDe.Range("K3").Value = "2016-09-04 20:51:41"
